$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44361,502,3858,397,6082,862,2083,1000,28500,0,0,0,3674,179,860,44197),
    @(44362,500,3958,1404,6092,254,2125,1000,28500,0,0,1232,3787,121,858,44462),
    @(44363,780,3808,1714,6752,74,2124,1000,28500,0,0,0,3787,76,857,44971),
    @(44364,1051,3883,1141,7093,413,2099,0,28500,0,0,0,3787,443,880,45361),
    @(44365,721,3554,1895,7756,535,2138,1500,28500,0,0,0,3787,0,819,45735),
    @(44368,527,3579,774,7625,768,2044,1000,28500,0,0,0,3787,170,809,45535),
    @(44369,200,3279,399,7574,204,1994,1500,28500,0,0,0,3787,133,821,45134),
    @(44370,705,3204,1325,7773,31,1951,1000,28500,0,0,0,3787,115,861,45215),
    @(44371,900,3053,1598,8132,299,1838,0,28500,0,0,0,3787,353,771,45309),
    @(44372,725,3057,1010,8331,527,1829,0,28500,0,0,0,3787,0,771,45503),
    @(44375,495,3025,553,8028,751,1812,1500,28500,0,0,0,3787,204,805,45152),
    @(44376,550,3375,479,7338,203,1811,1000,28500,0,0,0,3787,132,805,44811),
    @(44377,600,3270,1124,7338,15,1794,0,28500,0,0,0,3787,115,804,44689),
    @(44378,700,3070,1122,7398,281,1776,1000,28500,0,0,0,3787,281,732,44531),
    @(44379,640,2985,1156,7721,524,1774,1000,28500,0,0,0,3787,77,809,44766),
    @(44382,0,2985,237,7719,0,1774,0,28500,0,0,0,3787,0,809,44764),
    @(44383,800,3235,1275,7476,205,1775,1000,28500,0,0,0,3787,134,810,44773),
    @(44384,680,3315,1376,7749,14,1775,1000,28500,0,0,0,3787,45,740,45125),
    @(44385,1010,3625,1087,7177,288,1782,0,28500,0,0,0,3787,248,707,44871),
    @(44386,520,3505,587,6995,463,1721,0,28500,0,0,0,3787,0,630,44508),
    @(44389,237,3247,1029,6941,754,1724,1000,28500,0,0,0,3787,148,574,44199),
    @(44390,975,3422,831,7050,206,1726,1000,28500,0,0,0,3787,74,515,44485),
    @(44391,850,3592,987,6721,14,1726,0,28500,0,0,0,3787,43,514,44325),
    @(44393,1685,3747,1404,7000,763,1738,1000,28500,0,0,0,3787,252,517,44771),
    @(44396,0,3747,0,7000,0,1738,0,28500,0,0,0,3787,0,517,44771),
    @(44403,612,4122,1096,6688,768,1752,1000,28500,0,0,0,3787,148,517,44849),
    @(44404,1020,4167,997,6957,200,1745,1000,28500,0,0,0,3787,76,520,45156)
)

$startRow = 113
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 1; $c -le 16; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c-1]
    }
}
